$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel does not convert them to numbers
$textCells = @("D5","D6","D10","D11","D14","D17","D20","D21","D22","D24","D25","D27","D30","D33","D37","D38","D39","D40","D42","D44","D45","D46","D47","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '60.679.25'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '2.904.24'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '585.94'
$ws.Range("D6").Value = '147.30'
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D9").Value = '2.902.94'
$ws.Range("E9").Value = '  -3.75%  '
$ws.Range("D10").Value = '6.70'
$ws.Range("E10").Value = '  +4.55%  '
$ws.Range("D11").Value = '0.143'
$ws.Range("E11").Value = '  -4.23%  '
$ws.Range("E12").Value = '  -2.57%  '
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '34.01'
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").Value = '3.386.32'
$ws.Range("D17").Value = '6.81'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("D18").Value = '60.568.57'
$ws.Range("E18").Value = '  -2.71%  '
$ws.Range("D19").Value = '2.903.20'
$ws.Range("E19").Value = '  -3.83%  '
$ws.Range("D20").Value = '427.67'
$ws.Range("E20").Value = '  -4.73%  '
$ws.Range("D21").Value = '13.61'
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("D22").Value = '0.670'
$ws.Range("E22").Value = '  -3.00%  '
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("D24").Value = '80.68'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("D25").Value = '11.06'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").Value = '11.86'
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  +2.14%  '
$ws.Range("E32").Value = '  -3.27%  '
$ws.Range("D33").Value = '26.46'
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").Value = '0.0₃0837'
$ws.Range("E35").Value = '  -1.77%  '
$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("D37").Value = '5.67'
$ws.Range("E37").Value = '  -3.13%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '2.03'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '49.29'
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '2.96'
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").Value = '8.73'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("D44").Value = '41.81'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").Value = '0.0347'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '370.08'
$ws.Range("E46").Value = '  -6.11%  '
$ws.Range("D47").Value = '133.77'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").Value = '2.654.39'
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("D50").Value = '24.87'
$ws.Range("E50").Value = '  +4.75%  '
$ws.Range("E51").Value = '  -1.27%  '

# Restore default style on the forced-text cells (keep number format but drop explicit style index)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
